# Append " (deleted)" after the trailing "position" in the
# "view_theme(position, vars): View HTML of position" bullet, with the
# word "deleted" rendered in italic + light-gray highlight (the "delete"
# and "d" portions match the two runs introduced by the diff).

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("View HTML of position", $true, $false, $false, `
                            $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target text 'View HTML of position'"
}

# Collapse the found range to its end (right after the trailing "position").
$rng.Collapse(0)   # wdCollapseEnd

# Plain run: " ("
$rng.InsertAfter(" (")
$rng.Collapse(0)

# Italic + light-gray-highlighted run: "delete"
$rng.InsertAfter("delete")
$rng.Font.Italic = $true
$rng.Font.HighlightColorIndex = 16   # wdGray25 -> OOXML w:highlight="lightGray"
$rng.Collapse(0)

# Italic + light-gray-highlighted run: "d"
$rng.InsertAfter("d")
$rng.Font.Italic = $true
$rng.Font.HighlightColorIndex = 16   # wdGray25 -> OOXML w:highlight="lightGray"
$rng.Collapse(0)

# Plain run: ")"
$rng.InsertAfter(")")
